$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts LongName..Reject_Reason right by one)
$ws.Columns.Item(2).Insert()

# New column B = ParentId
$ws.Cells.Item(1, 2).Value = "ParentId"
$ws.Range("B2:B7").Value = 0

# Remove the trailing Reject_Reason column (now shifted to column AA)
$ws.Columns.Item(27).Delete()

# Update the selection to match the new edit location
$ws.Range("B1").Select()
